$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BudgetaForm")
$ws.Activate()
Write-Host ("A1: " + $ws.Range("A1").Value)
Write-Host ("A10: " + $ws.Range("A10").Value)
